$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.485.08"
$ws.Range("E2").Value = "  -4.36%  "
$ws.Range("D3").Value = "2.534.90"
$ws.Range("E3").Value = "  -4.07%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'507.10"
$ws.Range("E5").Value = "  -4.82%  "
$ws.Range("D6").Value = "'143.96"
$ws.Range("E6").Value = "  -8.03%  "
$ws.Range("E7").Value = "  +0.17%  "
$ws.Range("D8").Value = "'0.564"
$ws.Range("D9").Value = "2.539.47"
$ws.Range("E9").Value = "  -4.36%  "
$ws.Range("D10").Value = "'6.18"
$ws.Range("E10").Value = "  -7.63%  "
$ws.Range("E11").Value = "  -7.74%  "
$ws.Range("E12").Value = "  -5.58%  "
$ws.Range("E13").Value = "  -0.75%  "
$ws.Range("D14").Value = "2.982.56"
$ws.Range("E14").Value = "  -3.85%  "
$ws.Range("D15").Value = "58.474.48"
$ws.Range("E15").Value = "  -4.36%  "
$ws.Range("D16").Value = "'20.70"
$ws.Range("E16").Value = "  -6.16%  "
$ws.Range("E17").Value = "  -7.19%  "
$ws.Range("D18").Value = "2.544.63"
$ws.Range("E18").Value = "  -3.66%  "
$ws.Range("D19").Value = "'4.53"
$ws.Range("E19").Value = "  -5.67%  "
$ws.Range("D20").Value = "'334.78"
$ws.Range("E20").Value = "  -6.29%  "
$ws.Range("D21").Value = "'10.07"
$ws.Range("E21").Value = "  -5.98%  "
$ws.Range("E22").Value = "  -0.08%  "
$ws.Range("D23").Value = "'5.94"
$ws.Range("E23").Value = "  -5.26%  "
$ws.Range("D24").Value = "'60.11"
$ws.Range("E24").Value = "  -2.93%  "
$ws.Range("E25").Value = "  -5.63%  "
$ws.Range("D26").Value = "'0.998"
$ws.Range("E26").Value = "  -0.18%  "
$ws.Range("E27").Value = "  -5.84%  "
$ws.Range("D28").Value = "2.651.27"
$ws.Range("E28").Value = "  -3.59%  "
$ws.Range("D29").Value = "0.0₃0784"
$ws.Range("E29").Value = "  -10.47%  "
$ws.Range("D30").Value = "'6.96"
$ws.Range("E30").Value = "  -6.80%  "
$ws.Range("E31").Value = "  +0.04%  "
$ws.Range("D32").Value = "'149.37"
$ws.Range("E32").Value = "  -1.28%  "
$ws.Range("D33").Value = "'5.84"
$ws.Range("E33").Value = "  -5.46%  "
$ws.Range("D34").Value = "'18.51"
$ws.Range("E34").Value = "  -5.60%  "
$ws.Range("D35").Value = "'1.53"
$ws.Range("E35").Value = "  -5.88%  "
$ws.Range("D36").Value = "'0.927"
$ws.Range("E36").Value = "  +4.47%  "
$ws.Range("D37").Value = "'3.89"
$ws.Range("E38").Value = "  -8.13%  "
$ws.Range("D39").Value = "'36.01"
$ws.Range("E39").Value = "  -1.89%  "
$ws.Range("D40").Value = "'0.824"
$ws.Range("E40").Value = "  -11.43%  "
$ws.Range("E41").Value = "  -7.35%  "
$ws.Range("D42").Value = "'283.10"
$ws.Range("E42").Value = "  -5.04%  "
$ws.Range("E43").Value = "  -8.10%  "
$ws.Range("D44").Value = "'0.0995"
$ws.Range("E44").Value = "  -2.86%  "
$ws.Range("E45").Value = "  +0.06%  "
$ws.Range("D46").Value = "'0.601"
$ws.Range("E46").Value = "  -6.96%  "
$ws.Range("D47").Value = "'0.0533"
$ws.Range("E47").Value = "  -5.79%  "
$ws.Range("E48").Value = "  -6.37%  "
$ws.Range("D49").Value = "'10.30"
$ws.Range("E49").Value = "  -0.61%  "
$ws.Range("E50").Value = "  -5.70%  "
$ws.Range("E51").Value = "  -11.79%  "
